$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style of A16 (bordered/bold/centered) down to the new rows A17:A19
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write header/data rows 10-19 with final values
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.166259248206962
$ws.Cells.Item(10, 4).Value = 0.6150372372946338
$ws.Cells.Item(10, 5).Value = 1.04570840306212
$ws.Cells.Item(10, 6).Value = 1.166259248206962
$ws.Cells.Item(10, 7).Value = 0.7983980656971881
$ws.Cells.Item(10, 8).Value = 1.126180299928017
$ws.Cells.Item(10, 9).Value = 1.086610671014488
$ws.Cells.Item(10, 10).Value = 0.6150372372946338
$ws.Cells.Item(10, 11).Value = 0.8303728201783769
$ws.Cells.Item(10, 12).Value = 0.9983160341926693
$ws.Cells.Item(10, 13).Value = 0.9730323208672349

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 0.9995207151137521
$ws.Cells.Item(11, 4).Value = 0.9637772103659213
$ws.Cells.Item(11, 5).Value = 0.9822669232148015
$ws.Cells.Item(11, 6).Value = 0.9995207151137521
$ws.Cells.Item(11, 7).Value = 1.034150425731354
$ws.Cells.Item(11, 8).Value = 0.8694875128817304
$ws.Cells.Item(11, 9).Value = 0.9857289973144043
$ws.Cells.Item(11, 10).Value = 0.9637772103659213
$ws.Cells.Item(11, 11).Value = 0.9730220667903614
$ws.Cells.Item(11, 12).Value = 0.9862713909520567
$ws.Cells.Item(11, 13).Value = 0.9724886307703273

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 0.9972166131034629
$ws.Cells.Item(12, 4).Value = 0.9654221358289631
$ws.Cells.Item(12, 5).Value = 0.9827416855265932
$ws.Cells.Item(12, 6).Value = 0.9972166131034629
$ws.Cells.Item(12, 7).Value = 1.035113803753338
$ws.Cells.Item(12, 8).Value = 0.8697235836166758
$ws.Cells.Item(12, 9).Value = 0.9852586586559313
$ws.Cells.Item(12, 10).Value = 0.9654221358289631
$ws.Cells.Item(12, 11).Value = 0.9740819106777782
$ws.Cells.Item(12, 12).Value = 0.9856492618906205
$ws.Cells.Item(12, 13).Value = 0.9725794134141608

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 0.9993258281154973
$ws.Cells.Item(13, 4).Value = 0.9637924972707507
$ws.Cells.Item(13, 5).Value = 0.9824539607080326
$ws.Cells.Item(13, 6).Value = 0.9993258281154973
$ws.Cells.Item(13, 7).Value = 1.034244667105257
$ws.Cells.Item(13, 8).Value = 0.8692819243625699
$ws.Cells.Item(13, 9).Value = 0.985569647148516
$ws.Cells.Item(13, 10).Value = 0.9637924972707507
$ws.Cells.Item(13, 11).Value = 0.9731232289893916
$ws.Cells.Item(13, 12).Value = 0.9862245285524445
$ws.Cells.Item(13, 13).Value = 0.9724447541184372

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 0.7665679999999983
$ws.Cells.Item(14, 4).Value = 0.540900000000001
$ws.Cells.Item(14, 5).Value = 1.259892
$ws.Cells.Item(14, 6).Value = 0.7665679999999983
$ws.Cells.Item(14, 7).Value = 0.6358840000000004
$ws.Cells.Item(14, 8).Value = 1.933344000000001
$ws.Cells.Item(14, 9).Value = 1.162999999999998
$ws.Cells.Item(14, 10).Value = 0.540900000000001
$ws.Cells.Item(14, 11).Value = 0.9003960000000003
$ws.Cells.Item(14, 12).Value = 0.8334819999999993
$ws.Cells.Item(14, 13).Value = 1.049931333333333

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 0.72
$ws.Cells.Item(15, 4).Value = 0.07000000000000001
$ws.Cells.Item(15, 5).Value = 1.461237500000003
$ws.Cells.Item(15, 6).Value = 0.72
$ws.Cells.Item(15, 7).Value = 0.25
$ws.Cells.Item(15, 8).Value = 2.768750000000002
$ws.Cells.Item(15, 9).Value = 1.34
$ws.Cells.Item(15, 10).Value = 0.07000000000000001
$ws.Cells.Item(15, 11).Value = 0.7656187500000017
$ws.Cells.Item(15, 12).Value = 0.7428093750000009
$ws.Cells.Item(15, 13).Value = 1.101664583333334

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 0.8437880829951994
$ws.Cells.Item(16, 4).Value = 0.4525840355328005
$ws.Cells.Item(16, 5).Value = 1.266845427199999
$ws.Cells.Item(16, 6).Value = 0.8437880829951994
$ws.Cells.Item(16, 7).Value = 0.5616876208128014
$ws.Cells.Item(16, 8).Value = 2.008984824831995
$ws.Cells.Item(16, 9).Value = 1.191338608435194
$ws.Cells.Item(16, 10).Value = 0.4525840355328005
$ws.Cells.Item(16, 11).Value = 0.8597147313664
$ws.Cells.Item(16, 12).Value = 0.8517514071807997
$ws.Cells.Item(16, 13).Value = 1.054204766634665

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9950240908230863
$ws.Cells.Item(17, 4).Value = 0.9945358124397078
$ws.Cells.Item(17, 5).Value = 0.9939644677708043
$ws.Cells.Item(17, 6).Value = 0.9950240908230863
$ws.Cells.Item(17, 7).Value = 0.990812887787386
$ws.Cells.Item(17, 8).Value = 0.9945238266094194
$ws.Cells.Item(17, 9).Value = 0.9944208768167794
$ws.Cells.Item(17, 10).Value = 0.9945358124397078
$ws.Cells.Item(17, 11).Value = 0.994250140105256
$ws.Cells.Item(17, 12).Value = 0.9946371154641712
$ws.Cells.Item(17, 13).Value = 0.9938803270411971

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9219342028667733
$ws.Cells.Item(18, 4).Value = 1.055510848183846
$ws.Cells.Item(18, 5).Value = 1.005691515749741
$ws.Cells.Item(18, 6).Value = 0.9219342028667733
$ws.Cells.Item(18, 7).Value = 1.019793203122227
$ws.Cells.Item(18, 8).Value = 1.015292691945201
$ws.Cells.Item(18, 9).Value = 0.9830502005883052
$ws.Cells.Item(18, 10).Value = 1.055510848183846
$ws.Cells.Item(18, 11).Value = 1.030601181966794
$ws.Cells.Item(18, 12).Value = 0.9762676924167835
$ws.Cells.Item(18, 13).Value = 1.000212110409349

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9721686643051966
$ws.Cells.Item(19, 4).Value = 1.155034933462931
$ws.Cells.Item(19, 5).Value = 0.952136798475879
$ws.Cells.Item(19, 6).Value = 0.9721686643051966
$ws.Cells.Item(19, 7).Value = 1.080038061395584
$ws.Cells.Item(19, 8).Value = 0.8775173774800933
$ws.Cells.Item(19, 9).Value = 0.9524183233782294
$ws.Cells.Item(19, 10).Value = 1.155034933462931
$ws.Cells.Item(19, 11).Value = 1.053585865969405
$ws.Cells.Item(19, 12).Value = 1.012877265137301
$ws.Cells.Item(19, 13).Value = 0.998219026416319
